# Auto-generated script applying numeric corrections to the Golem_Profits
# (per-sheet leve-profit) tables, as scraped by the scheduled market-data runner.
# Each block updates one leve row's price/profit columns (H:N) to reflect
# refreshed Market Board averages; a handful of rows also gain/lose the
# HQ-profit (N) or NQ-profit (M) cell depending on whether an HQ price exists.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 83.5
$ws.Range("I38").Value = 83.5
$ws.Range("K38").Value = 250.5
$ws.Range("M38").Value = 121.5

$ws.Range("H39").Value = 126
$ws.Range("I39").Value = 91.2
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 273.6
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = 22.39999999999998
$ws.Range("N39").Value = -1492

$ws.Range("H58").Value = 608.1667
$ws.Range("I58").Value = 608.1667
$ws.Range("K58").Value = 1824.5001
$ws.Range("M58").Value = -1674.5001

$ws.Range("H92").Value = 271.07693
$ws.Range("I92").Value = 261.27274
$ws.Range("K92").Value = 261.27274
$ws.Range("M92").Value = 986.72726

$ws.Range("H97").Value = 748.8889
$ws.Range("J97").Value = 748.8889
$ws.Range("L97").Value = 2246.6667
$ws.Range("N97").Value = -3238.6667

$ws.Range("H100").Value = 2472
$ws.Range("I100").Value = 1444
$ws.Range("K100").Value = 1444
$ws.Range("M100").Value = -903

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 59
$ws.Range("J5").Value = 15
$ws.Range("L5").Value = 15
$ws.Range("N5").Value = -239

$ws.Range("H97").Value = 2810
$ws.Range("I97").Value = 2500
$ws.Range("K97").Value = 2500
$ws.Range("M97").Value = -2004

$ws.Range("H135").Value = 5024984.5
$ws.Range("J135").Value = 5024984.5
$ws.Range("L135").Value = 5024984.5
$ws.Range("N135").Value = -5035124.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 59
$ws.Range("J4").Value = 15
$ws.Range("L4").Value = 15
$ws.Range("N4").Value = -245

$ws.Range("H9").Value = 49999
$ws.Range("J9").Value = 49999
$ws.Range("L9").Value = 49999
$ws.Range("N9").Value = -50335

$ws.Range("H20").Value = 1034.4
$ws.Range("I20").Value = 1034.4
$ws.Range("K20").Value = 1034.4
$ws.Range("M20").Value = -787.4000000000001

$ws.Range("H94").Value = 2366.5334
$ws.Range("I94").Value = 1874.875
$ws.Range("K94").Value = 1874.875
$ws.Range("M94").Value = -1423.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1199.6666
$ws.Range("I16").Value = 1199.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1199.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -912.6666
$ws.Range("N16").ClearContents()

$ws.Range("H113").Value = 1199.6666
$ws.Range("I113").Value = 1199.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1199.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 970.3334
$ws.Range("N113").ClearContents()

$ws.Range("H134").Value = 1373.1428
$ws.Range("I134").Value = 1373.1428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4119.428400000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1584.428400000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5600
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5600
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 16800
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -17024

$ws.Range("H7").Value = 184.75
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 196.85715
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 590.5714499999999
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -814.5714499999999

$ws.Range("H86").Value = 1753.8
$ws.Range("J86").Value = 1753.8
$ws.Range("L86").Value = 5261.4
$ws.Range("N86").Value = -7633.4

$ws.Range("H89").Value = 1753.8
$ws.Range("J89").Value = 1753.8
$ws.Range("L89").Value = 15784.2
$ws.Range("N89").Value = -27640.2

$ws.Range("H133").Value = 516.5
$ws.Range("J133").Value = 516.5
$ws.Range("L133").Value = 1549.5
$ws.Range("N133").Value = -11669.5

$ws.Range("H138").Value = 1450
$ws.Range("I138").Value = 1450
$ws.Range("K138").Value = 4350
$ws.Range("M138").Value = 790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 41000
$ws.Range("J15").Value = 41000
$ws.Range("L15").Value = 41000
$ws.Range("N15").Value = -41576

$ws.Range("H48").Value = 24997
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 24997
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 24997
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -25967

$ws.Range("H80").Value = 2400
$ws.Range("J80").Value = 2400
$ws.Range("L80").Value = 2400
$ws.Range("N80").Value = -4396

$ws.Range("H81").Value = 41000
$ws.Range("J81").Value = 41000
$ws.Range("L81").Value = 41000
$ws.Range("N81").Value = -42996

$ws.Range("H83").Value = 2400
$ws.Range("J83").Value = 2400
$ws.Range("L83").Value = 12000
$ws.Range("N83").Value = -21984

$ws.Range("H84").Value = 41000
$ws.Range("J84").Value = 41000
$ws.Range("L84").Value = 123000
$ws.Range("N84").Value = -132984

$ws.Range("H97").Value = 304.5
$ws.Range("I97").Value = 109
$ws.Range("K97").Value = 109
$ws.Range("M97").Value = 387

$ws.Range("H113").Value = 3169.2727
$ws.Range("I113").Value = 2484.875
$ws.Range("K113").Value = 2484.875
$ws.Range("M113").Value = -314.875

$ws.Range("H132").Value = 1071.8
$ws.Range("I132").Value = 743.5
$ws.Range("K132").Value = 2230.5
$ws.Range("M132").Value = 299.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5068.7
$ws.Range("I7").Value = 4928.2856
$ws.Range("K7").Value = 4928.2856
$ws.Range("M7").Value = -4816.2856

$ws.Range("H61").Value = 2709.5
$ws.Range("I61").Value = 2888.3333
$ws.Range("K61").Value = 2888.3333
$ws.Range("M61").Value = -2686.3333

$ws.Range("H82").Value = 1409.5454
$ws.Range("I82").Value = 1350.5
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1350.5
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -989.5
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 1409.5454
$ws.Range("I85").Value = 1350.5
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1350.5
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -102.5
$ws.Range("N85").Value = -4496

$ws.Range("H113").Value = 2709.5
$ws.Range("I113").Value = 2888.3333
$ws.Range("K113").Value = 2888.3333
$ws.Range("M113").Value = -718.3332999999998

$ws.Range("H126").Value = 5068.7
$ws.Range("I126").Value = 4928.2856
$ws.Range("K126").Value = 14784.8568
$ws.Range("M126").Value = -12314.8568

$ws.Range("H132").Value = 2253.8
$ws.Range("I132").Value = 1942.25
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5826.75
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3296.75
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2527

$ws.Range("H48").Value = 10059
$ws.Range("I48").Value = 10059
$ws.Range("K48").Value = 10059
$ws.Range("M48").Value = -9490

$ws.Range("H113").Value = 4250.5
$ws.Range("I113").Value = 338
$ws.Range("K113").Value = 1014
$ws.Range("M113").Value = 1156
